$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row describing the "delete a given (non-tail) linked-list node" problem
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 237
$ws.Cells.Item(16, 3).Value = "请编写一个函数，使其可以删除某个链表中给定的（非末尾）节点，你将只被给定要求被删除的节点。 "
$ws.Cells.Item(16, 4).Value = "1 目标节点cur的下一个节点nextNode`n2 nextNode的值赋值给cur节点`n3 删除nextNode节点即可`n4 如果cur是尾节点，将cur=null【待确认】"
$ws.Cells.Item(16, 5).Value = "删除节点"
$ws.Cells.Item(16, 6).Value = "O(1)"
$ws.Cells.Item(16, 7).Value = "O(1)"

# Match the formatting used by the other data rows (Calibri 14, left/top align, wrap text)
$rng = $ws.Range("A16:G16")
$rng.Font.Name = "Calibri"
$rng.Font.Size = 14
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160
$rng.WrapText = $true

$ws.Rows.Item(16).RowHeight = 80

$ws.Range("D18").Select()
